$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------
# 1. Duplicate the existing "fiche" template (rows 20:35) twice to create
#    the skeletons for the two new product sheets ("Controle parental
#    Swisscom" at rows 40:58 and "Internet Security Swisscom" at rows 61:79)
# ----------------------------------------------------------------------
$ws.Range("B20:C35").Copy($ws.Range("B40"))
$excel.CutCopyMode = 0

# Insert the 3 extra "bordered" data rows used by the new fiche (rows 54:56)
# replicating the border/format carried by the last template data row (53)
$ws.Rows.Item(54).Insert()
$ws.Rows.Item(54).Insert()
$ws.Rows.Item(54).Insert()
$ws.Range("B53:C53").Copy()
$ws.Range("B54:C56").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Second fiche skeleton, pasted after a blank separator row (row 60)
$ws.Range("B20:C35").Copy($ws.Range("B61"))
$excel.CutCopyMode = 0
$ws.Rows.Item(75).Insert()
$ws.Rows.Item(75).Insert()
$ws.Rows.Item(75).Insert()
$ws.Range("B74:C74").Copy()
$ws.Range("B75:C77").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ----------------------------------------------------------------------
# 2. Fill in the text content for "Controle parental Swisscom" (rows 40:58)
# ----------------------------------------------------------------------
$ws.Range("B40").Value2 = "Contrôle parental Swisscom"
$ws.Range("B41").Value2 = "Fonctionnement"
$ws.Range("B42").Value2 = "L'utilisateur choisit les appareils sur lesquels il veut activer le contrôle parental. Celui-ci consiste dans cette version standard à limiter les heures de navigation sur internet individuellement. Un créneau horaire peut être défini selon les jours de la semaine et du weekend ainsi que le temps total. Ce contrôle est paramétrable sur la swisscom TV box. Les émissions contenant des limites d'âge peuvent également être bloquées par un PIN que seules les personnes autorisées connaissent."
$ws.Range("B43").Value2 = "Critère"
$ws.Range("C43").Value2 = "Observation"
$ws.Range("B44").Value2 = "Blocage"
$ws.Range("C44").Value2 = "Bloque la navigation sur internet et certains contenus"
$ws.Range("B45").Value2 = "Orienté web/natif/mobile"
$ws.Range("C45").Value2 = "Orienté web"
$ws.Range("B46").Value2 = "Plateforme"
$ws.Range("C46").Value2 = "Fonctionne pour tout ordinateur, tablette et smartphone pouvant se connecter à internet"
$ws.Range("B47").Value2 = "Base volontaire ou restrictions"
$ws.Range("C47").Value2 = "Ne fonctionne que sur la base de restrictions. L'utilisateur du réseau contrôlé ne peut pas changer quoi que ce soit aux restrictions, il subit des interdits. "
$ws.Range("B48").Value2 = "Public cible"
$ws.Range("C48").Value2 = "Les enfants particulièrement`net les adolescents"
$ws.Range("B49").Value2 = "Paramétrer en fonction de l'âge"
$ws.Range("C49").Value2 = "Navigation internet -> Pas directement. Si l'enfant utilise toujours les mêmes appareils, alors le contrôle est indirectement associé à l'âge.`nBlocage de contenu lié à l'âge -> Oui"
$ws.Range("B50").Value2 = "Redonne le contrôle"
$ws.Range("C50").Value2 = "Ne donne aucun contrôle. Le service`ninterdit mais n'apprend pas à gérer "
$ws.Range("B51").Value2 = "Difficulté de prise en main"
$ws.Range("C51").Value2 = "Très facile à utiliser"
$ws.Range("B52").Value2 = "Administration depuis la même machine"
$ws.Range("C52").Value2 = "La gateway est paramétrée via une interface web donc administrable depuis n'importe quelle machine"
$ws.Range("B53").Value2 = "Fournit des statistiques d'utilisation"
$ws.Range("C53").Value2 = "Non, il ne fournit pas de statistiques"
$ws.Range("B54").Value2 = "Avis des utilisateurs"
$ws.Range("C54").Value2 = "Aucun avis disponible"
$ws.Range("B55").Value2 = "Payant"
$ws.Range("C55").Value2 = "La fonctionnalité est comprise dans l'abonnement swisscom"
$ws.Range("B56").Value2 = "Offline"
$ws.Range("C56").Value2 = "Est inutile en hors ligne puisqu’il monitore`nuniquement ce qui est en ligne"
$ws.Range("B57").Value2 = "Remarque complémentaire"
$ws.Range("B58").Value2 = "Ce contrôle parental est facile d'utilisation. Par contre, il ne permet pas de filtrer spécifiquement le contenu (mis à part le blocage du contenu limité par l'âge). Internet security permet un contrôle beaucoup plus approfondi, c'est le prochain service analysé."

# ----------------------------------------------------------------------
# 3. Fill in the text content for "Internet Security Swisscom" (rows 61:79)
# ----------------------------------------------------------------------
$ws.Range("B61").Value2 = "Internet Security Swisscom"
$ws.Range("B62").Value2 = "Fonctionnement"
$ws.Range("B63").Value2 = "L'utilisateur installe cette application sur la machine que l'enfant ou l'adolescent utilise depuis le compte administrateur. L'application définit le contenu non-accessible par catégories et par sites spécifiques. Elle permet aussi de définir les horaires de navigation sur internet et d'accès à la machine même."
$ws.Range("B64").Value2 = "Critère"
$ws.Range("C64").Value2 = "Observation"
$ws.Range("B65").Value2 = "Blocage"
$ws.Range("C65").Value2 = "Bloque l'accès à l'ordinateur, les heures de navigation sur internet et certains contenus"
$ws.Range("B66").Value2 = "Orienté web/natif/mobile"
$ws.Range("C66").Value2 = "Application native"
$ws.Range("B67").Value2 = "Plateforme"
$ws.Range("C67").Value2 = "Fonctionne sur ordinateur, tablette et smartphone (fonctionnalités différentes)"
$ws.Range("B68").Value2 = "Base volontaire ou restrictions"
$ws.Range("C68").Value2 = "Ne fonctionne que sur la base de restrictions. L'utilisateur du réseau contrôlé ne peut pas changer quoi que ce soit aux restrictions, il subit des interdits. "
$ws.Range("B69").Value2 = "Public cible"
$ws.Range("C69").Value2 = "Les enfants particulièrement`net les adolescents"
$ws.Range("B70").Value2 = "Paramétrer en fonction de l'âge"
$ws.Range("C70").Value2 = "Non"
$ws.Range("B71").Value2 = "Redonne le contrôle"
$ws.Range("C71").Value2 = "Ne donne aucun contrôle. L'application`ninterdit mais n'apprend pas à gérer "
$ws.Range("B72").Value2 = "Difficulté de prise en main"
$ws.Range("C72").Value2 = "Facile à utiliser"
$ws.Range("B73").Value2 = "Administration depuis la même machine"
$ws.Range("C73").Value2 = "Les restrictions sont élaborées depuis la session administrateur de la machine"
$ws.Range("B74").Value2 = "Fournit des statistiques d'utilisation"
$ws.Range("C74").Value2 = "Non, il ne fournit pas de statistiques"
$ws.Range("B75").Value2 = "Avis des utilisateurs"
$ws.Range("C75").Value2 = "Bonne application (4.4/5 sur google play)"
$ws.Range("B76").Value2 = "Payant"
$ws.Range("C76").Value2 = "La fonctionnalité est gratuite pendant 6 mois et ensuite payante"
$ws.Range("B77").Value2 = "Offline"
$ws.Range("C77").Value2 = "Utile car restreint également l'accès à la machine et pas seulement internet."
$ws.Range("B78").Value2 = "Remarque complémentaire"
$ws.Range("B79").Value2 = "Cette application est également très utile pour la protection des données et contre les sites webs dangereux ainsi que les virus. "

# ----------------------------------------------------------------------
# 4. Row heights for the new fiches (and the blank separator rows)
# ----------------------------------------------------------------------
$ws.Rows.Item(39).RowHeight = 16.5
$ws.Rows.Item(41).RowHeight = 16.5
$ws.Rows.Item(42).RowHeight = 81
$ws.Rows.Item(44).RowHeight = 25.5
$ws.Rows.Item(46).RowHeight = 35.25
$ws.Rows.Item(47).RowHeight = 51
$ws.Rows.Item(48).RowHeight = 25.5
$ws.Rows.Item(49).RowHeight = 63.75
$ws.Rows.Item(50).RowHeight = 25.5
$ws.Rows.Item(52).RowHeight = 38.25
$ws.Rows.Item(55).RowHeight = 25.5
$ws.Rows.Item(56).RowHeight = 25.5
$ws.Rows.Item(58).RowHeight = 43.5
$ws.Rows.Item(60).RowHeight = 16.5
$ws.Rows.Item(62).RowHeight = 16.5
$ws.Rows.Item(63).RowHeight = 80.25
$ws.Rows.Item(65).RowHeight = 25.5
$ws.Rows.Item(67).RowHeight = 25.5
$ws.Rows.Item(68).RowHeight = 51
$ws.Rows.Item(69).RowHeight = 25.5
$ws.Rows.Item(71).RowHeight = 25.5
$ws.Rows.Item(73).RowHeight = 25.5
$ws.Rows.Item(76).RowHeight = 25.5
$ws.Rows.Item(77).RowHeight = 25.5
$ws.Rows.Item(79).RowHeight = 49.5

# ----------------------------------------------------------------------
# 5. Workbook / window level cosmetic updates
# ----------------------------------------------------------------------
$ws.Application.ActiveWindow.Zoom = 85
$ws.Range("F66").Select()

